# Generate Report for Handoff
# Adds two new file entries (558fb5c0-... and 82c96c97-...) as rows 4 & 5
# on all three sheets (Overview, zh-cn, de-de), and refreshes a couple of
# existing timestamp strings.

$wb = $excel.ActiveWorkbook

$mdBase   = "https://github.com/OpenLocalizationTest/oltest/blob/745f75da9ff980264e36474f844e7f7f80b429c8/e2e"
$zhBase   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5cdee8086d7f08b96dded9e42b501a8f6226dd7b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high"
$deBase   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ac36f27b77e64b9f045706ad1406f7a1e9ea3109/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high"

$file1 = "558fb5c0-f923-49b7-adf6-dae2c85f31d4"
$file2 = "82c96c97-4831-4d5b-b3b5-b6e8ab6e3b8a"

$zhHash1 = "095179d80fd8234f14f93d97d70fa30b048acb82"
$zhHash2 = "a64e5e4494b54e998ba4b39455c31528b12e25e2"

$md1 = "$file1.md"
$md2 = "$file2.md"

$zhXlf1 = "$file1.$zhHash1.zh-cn.xlf"
$zhXlf2 = "$file2.$zhHash2.zh-cn.xlf"
$deXlf1 = "$file1.$zhHash1.de-de.xlf"
$deXlf2 = "$file2.$zhHash2.de-de.xlf"

# ---------------------------------------------------------------------
# Sheet "Overview" — four columns: File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "$mdBase/$md1", "", "", $md1) | Out-Null
$wsOverview.Range("B4").Value = "Ready for handoff"
$wsOverview.Range("C4").Value = "Ready for handoff"
$wsOverview.Range("D4").Value = "2016-12-13 06:12:13"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), "$mdBase/$md2", "", "", $md2) | Out-Null
$wsOverview.Range("B5").Value = "Ready for handoff"
$wsOverview.Range("C5").Value = "Ready for handoff"
$wsOverview.Range("D5").Value = "2016-12-13 06:12:13"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "$mdBase/$md1", "", "", $md1) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B4"), "$mdBase/$md1", "", "", ".md") | Out-Null
$wsZh.Range("C4").Value = "Ready for handoff"
$wsZh.Hyperlinks.Add($wsZh.Range("D4"), "$zhBase/$zhXlf1", "", "", $zhXlf1) | Out-Null
$wsZh.Range("E4").Value = "2016-03-13 06:12:10"
$wsZh.Range("H4").Value = "0001-01-01 00:00:00"
$wsZh.Range("I4").Value = "Include"

$wsZh.Hyperlinks.Add($wsZh.Range("A5"), "$mdBase/$md2", "", "", $md2) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B5"), "$mdBase/$md2", "", "", ".md") | Out-Null
$wsZh.Range("C5").Value = "Ready for handoff"
$wsZh.Hyperlinks.Add($wsZh.Range("D5"), "$zhBase/$zhXlf2", "", "", $zhXlf2) | Out-Null
$wsZh.Range("E5").Value = "2016-03-13 06:12:10"
$wsZh.Range("H5").Value = "0001-01-01 00:00:00"
$wsZh.Range("I5").Value = "Include"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "$mdBase/$md1", "", "", $md1) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B4"), "$mdBase/$md1", "", "", ".md") | Out-Null
$wsDe.Range("C4").Value = "Ready for handoff"
$wsDe.Hyperlinks.Add($wsDe.Range("D4"), "$deBase/$deXlf1", "", "", $deXlf1) | Out-Null
$wsDe.Range("E4").Value = "2016-03-13 06:12:13"
$wsDe.Range("H4").Value = "0001-01-01 00:00:00"
$wsDe.Range("I4").Value = "Include"

$wsDe.Hyperlinks.Add($wsDe.Range("A5"), "$mdBase/$md2", "", "", $md2) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B5"), "$mdBase/$md2", "", "", ".md") | Out-Null
$wsDe.Range("C5").Value = "Ready for handoff"
$wsDe.Hyperlinks.Add($wsDe.Range("D5"), "$deBase/$deXlf2", "", "", $deXlf2) | Out-Null
$wsDe.Range("E5").Value = "2016-03-13 06:12:13"
$wsDe.Range("H5").Value = "0001-01-01 00:00:00"
$wsDe.Range("I5").Value = "Include"

# ---------------------------------------------------------------------
# Refresh the two timestamp strings that changed in the existing rows
# ---------------------------------------------------------------------
$wsOverview.Range("D2").Value = "2016-12-13 06:12:13"
$wsZh.Range("E2").Value = "2016-03-13 06:12:10"
$wsZh.Range("E3").Value = "2016-03-13 06:12:10"
$wsDe.Range("E2").Value = "2016-03-13 06:12:13"
$wsDe.Range("E3").Value = "2016-03-13 06:12:13"

"Report rows for $md1 and $md2 added to Overview, zh-cn and de-de sheets."
